# Updated cryptos list on Mon Sep 25 13:08:04 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.135.67"
$ws.Range("E2").Value = "  -2.30%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.572.56"
$ws.Range("E3").Value = "  -1.95%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.41%  "

# Row 5 - BNB (numeric-looking text -> keep as text with leading quote)
$ws.Range("D5").Value = "'208.24"

# Row 6 - XRP
$ws.Range("E6").Value = "  -3.12%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.40%  "

# Row 8 - Dogecoin
$ws.Range("E8").Value = "  -1.76%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.36%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = "  -0.96%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0844"

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.793.54"

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.580.57"
$ws.Range("E13").Value = "  -1.46%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.61%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.513"

# Row 16 - Litecoin
$ws.Range("E16").Value = "  -1.23%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.121.76"
$ws.Range("E17").Value = "  -2.21%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -2.48%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'7.28"
$ws.Range("E19").Value = "  +1.73%  "

# Row 20
$ws.Range("D20").Value = "'207.48"
$ws.Range("E20").Value = "  -1.35%  "

# Row 21
$ws.Range("E21").Value = "  -0.37%  "

# Row 22
$ws.Range("E22").Value = "  -1.58%  "

# Row 23
$ws.Range("E23").Value = "  -2.95%  "

# Row 24
$ws.Range("E24").Value = "  -2.89%  "

# Row 25
$ws.Range("D25").Value = "'143.31"
$ws.Range("E25").Value = "  -0.30%  "

# Row 26
$ws.Range("E26").Value = "  -0.27%  "

# Row 27
$ws.Range("E27").Value = "  -2.07%  "

# Row 29
$ws.Range("D29").Value = "'15.19"
$ws.Range("E29").Value = "  -1.21%  "

# Row 30
$ws.Range("E30").Value = "  -0.60%  "

# Row 31
$ws.Range("E31").Value = "  -1.65%  "

# Row 32
$ws.Range("E32").Value = "  -2.11%  "

# Row 33
$ws.Range("D33").Value = "'2.98"
$ws.Range("E33").Value = "  +0.22%  "

# Row 34
$ws.Range("D34").Value = "1.273.54"
$ws.Range("E34").Value = "  -1.49%  "

# Row 35
$ws.Range("E35").Value = "  +3.00%  "

# Row 36
$ws.Range("E36").Value = "  -1.50%  "

# Row 38
$ws.Range("E38").Value = "  -2.81%  "

# Row 39
$ws.Range("E39").Value = "  -10.15%  "

# Row 40
$ws.Range("D40").Value = "'0.809"
$ws.Range("E40").Value = "  -2.78%  "

# Row 41
$ws.Range("D41").Value = "'5.54"
$ws.Range("E41").Value = "  +1.71%  "

# Row 42
$ws.Range("E42").Value = "  -2.89%  "

# Row 43 / Row 44 swap: Aave <-> TrustWalletToken
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.762"
$ws.Range("E43").Value = "  -2.40%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'61.94"
$ws.Range("E44").Value = "  -2.04%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.707.02"
$ws.Range("E45").Value = "  -1.97%  "

# Row 46 - Quant
$ws.Range("D46").Value = "'88.94"
$ws.Range("E46").Value = "  -1.62%  "

# Rows 47/48/49 rotation: RenderToken -> BabyDogeCoin, Algorand -> RenderToken, BabyDogeCoin -> Algorand
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0104"
$ws.Range("E47").Value = "  +0.16%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.51"
$ws.Range("E48").Value = "  -3.40%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.100"
$ws.Range("E49").Value = "  -2.33%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -1.67%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  -0.30%  "
